# Update the "partial" sliding-window results sheet: recompute the
# IPC PO (col C), DELTA (col D) and DELTA^2 (col E) columns for every
# data row, a handful of IPC RO (col B) values whose float
# representation shifted slightly, and the TOTAL/MSE summary cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 29.81136131286621
$ws.Cells.Item(2, 4).Value = -0.02863868713379247
$ws.Cells.Item(2, 5).Value = 0.0008201744007472505
$ws.Cells.Item(3, 3).Value = 29.75001335144043
$ws.Cells.Item(3, 4).Value = -0.05998664855957259
$ws.Cells.Item(3, 5).Value = 0.003598398005409672
$ws.Cells.Item(4, 3).Value = 29.76519584655762
$ws.Cells.Item(4, 4).Value = -0.1548041534423845
$ws.Cells.Item(4, 5).Value = 0.02396432592301333
$ws.Cells.Item(5, 3).Value = 29.79987907409668
$ws.Cells.Item(5, 4).Value = -0.1801209259033243
$ws.Cells.Item(5, 5).Value = 0.03244354794827084
$ws.Cells.Item(6, 2).Value = 30.03999999999999
$ws.Cells.Item(6, 3).Value = 30.01015281677246
$ws.Cells.Item(6, 4).Value = -0.0298471832275311
$ws.Cells.Item(6, 5).Value = 0.000890854346617814
$ws.Cells.Item(7, 2).Value = 30.21000000000001
$ws.Cells.Item(7, 3).Value = 30.06319236755371
$ws.Cells.Item(7, 4).Value = -0.146807632446297
$ws.Cells.Item(7, 5).Value = 0.02155248094448704
$ws.Cells.Item(8, 3).Value = 30.18989372253418
$ws.Cells.Item(8, 4).Value = -0.03010627746581918
$ws.Cells.Item(8, 5).Value = 0.0009063879428488914
$ws.Cells.Item(9, 3).Value = 30.25397300720215
$ws.Cells.Item(9, 4).Value = -0.126026992797847
$ws.Cells.Item(9, 5).Value = 0.01588280291366858
$ws.Cells.Item(10, 3).Value = 30.50796318054199
$ws.Cells.Item(10, 4).Value = 0.06796318054199446
$ws.Cells.Item(10, 5).Value = 0.004618993909383734
$ws.Cells.Item(11, 3).Value = 30.37946128845215
$ws.Cells.Item(11, 4).Value = -0.1005387115478555
$ws.Cells.Item(11, 5).Value = 0.0101080325197029
$ws.Cells.Item(12, 3).Value = 30.38317680358887
$ws.Cells.Item(12, 4).Value = -0.3068231964111305
$ws.Cells.Item(12, 5).Value = 0.09414047385594319
$ws.Cells.Item(13, 3).Value = 30.47823143005371
$ws.Cells.Item(13, 4).Value = -0.2717685699462891
$ws.Cells.Item(13, 5).Value = 0.07385815561065101
$ws.Cells.Item(14, 3).Value = 30.61132621765137
$ws.Cells.Item(14, 4).Value = -0.3286737823486305
$ws.Cells.Item(14, 5).Value = 0.108026455203355
$ws.Cells.Item(15, 3).Value = 30.73988914489746
$ws.Cells.Item(15, 4).Value = -0.2101108551025419
$ws.Cells.Item(15, 5).Value = 0.04414657143192136
$ws.Cells.Item(16, 3).Value = 31.15546226501465
$ws.Cells.Item(16, 4).Value = 0.1354622650146524
$ws.Cells.Item(16, 5).Value = 0.01835002524289992
$ws.Cells.Item(17, 3).Value = 31.3719310760498
$ws.Cells.Item(17, 4).Value = 0.2519310760498001
$ws.Cells.Item(17, 5).Value = 0.06346926707961018
$ws.Cells.Item(18, 3).Value = 31.45668601989746
$ws.Cells.Item(18, 4).Value = 0.1766860198974598
$ws.Cells.Item(18, 5).Value = 0.03121794962720556
$ws.Cells.Item(19, 3).Value = 31.35834312438965
$ws.Cells.Item(19, 4).Value = -0.02165687561034702
$ws.Cells.Item(19, 5).Value = 0.0004690202612020434
$ws.Cells.Item(20, 3).Value = 31.55037498474121
$ws.Cells.Item(20, 4).Value = -0.02962501525878736
$ws.Cells.Item(20, 5).Value = 0.0008776415290833837
$ws.Cells.Item(21, 2).Value = 31.65000000000001
$ws.Cells.Item(21, 3).Value = 32.06996536254883
$ws.Cells.Item(21, 4).Value = 0.4199653625488224
$ws.Cells.Item(21, 5).Value = 0.1763709057407639
$ws.Cells.Item(22, 3).Value = 32.6776237487793
$ws.Cells.Item(22, 4).Value = 0.7976237487793014
$ws.Cells.Item(22, 5).Value = 0.6362036446167462
$ws.Cells.Item(23, 3).Value = 32.63322830200195
$ws.Cells.Item(23, 4).Value = 0.353228302001952
$ws.Cells.Item(23, 5).Value = 0.1247702333351822
$ws.Cells.Item(24, 3).Value = 32.70923233032227
$ws.Cells.Item(24, 4).Value = 0.2592323303222628
$ws.Cells.Item(24, 5).Value = 0.06720140108431076
$ws.Cells.Item(25, 2).Value = 32.84999999999999
$ws.Cells.Item(25, 3).Value = 32.97509384155273
$ws.Cells.Item(25, 4).Value = 0.1250938415527401
$ws.Cells.Item(25, 5).Value = 0.01564846919442204
$ws.Cells.Item(26, 2).Value = 32.90000000000001
$ws.Cells.Item(26, 3).Value = 33.11442184448242
$ws.Cells.Item(26, 4).Value = 0.2144218444824162
$ws.Cells.Item(26, 5).Value = 0.04597672739124147
$ws.Cells.Item(27, 2).Value = 33.09999999999999
$ws.Cells.Item(27, 3).Value = 33.00744247436523
$ws.Cells.Item(27, 4).Value = -0.09255752563475994
$ws.Cells.Item(27, 5).Value = 0.008566895551629244
$ws.Cells.Item(28, 2).Value = 33.40000000000001
$ws.Cells.Item(28, 3).Value = 33.61306381225586
$ws.Cells.Item(28, 4).Value = 0.2130638122558537
$ws.Cells.Item(28, 5).Value = 0.04539618809299767
$ws.Cells.Item(29, 3).Value = 33.67734146118164
$ws.Cells.Item(29, 4).Value = -0.02265853881836222
$ws.Cells.Item(29, 5).Value = 0.0005134093813832275
$ws.Cells.Item(30, 2).Value = 34.09999999999999
$ws.Cells.Item(30, 3).Value = 33.8823127746582
$ws.Cells.Item(30, 4).Value = -0.2176872253417912
$ws.Cells.Item(30, 5).Value = 0.04738772807700778
$ws.Cells.Item(31, 2).Value = 34.40000000000001
$ws.Cells.Item(31, 3).Value = 34.47198104858398
$ws.Cells.Item(31, 4).Value = 0.07198104858397869
$ws.Cells.Item(31, 5).Value = 0.0051812713552491
$ws.Cells.Item(32, 2).Value = 34.90000000000001
$ws.Cells.Item(32, 3).Value = 35.08790969848633
$ws.Cells.Item(32, 4).Value = 0.1879096984863224
$ws.Cells.Item(32, 5).Value = 0.03531005478522061
$ws.Cells.Item(33, 3).Value = 35.6751823425293
$ws.Cells.Item(33, 4).Value = 0.3751823425292997
$ws.Cells.Item(33, 5).Value = 0.1407617901457728
$ws.Cells.Item(34, 3).Value = 35.94461441040039
$ws.Cells.Item(34, 4).Value = 0.2446144104003878
$ws.Cells.Item(34, 5).Value = 0.05983620977552934
$ws.Cells.Item(35, 3).Value = 36.12253189086914
$ws.Cells.Item(35, 4).Value = -0.1774681091308565
$ws.Cells.Item(35, 5).Value = 0.0314949297584816
$ws.Cells.Item(36, 3).Value = 36.61664962768555
$ws.Cells.Item(36, 4).Value = -0.1833503723144503
$ws.Cells.Item(36, 5).Value = 0.03361735902784754
$ws.Cells.Item(37, 3).Value = 37.1579704284668
$ws.Cells.Item(37, 4).Value = -0.1420295715332003
$ws.Cells.Item(37, 5).Value = 0.02017239918990446
$ws.Cells.Item(38, 2).Value = 37.90000000000001
$ws.Cells.Item(38, 3).Value = 37.8542594909668
$ws.Cells.Item(38, 4).Value = -0.04574050903320881
$ws.Cells.Item(38, 5).Value = 0.002092194166617057
$ws.Cells.Item(39, 3).Value = 38.33866500854492
$ws.Cells.Item(39, 4).Value = -0.1613349914550781
$ws.Cells.Item(39, 5).Value = 0.02602897946781013
$ws.Cells.Item(40, 2).Value = 38.90000000000001
$ws.Cells.Item(40, 3).Value = 39.00360488891602
$ws.Cells.Item(40, 4).Value = 0.1036048889160099
$ws.Cells.Item(40, 5).Value = 0.01073397300729876
$ws.Cells.Item(41, 2).Value = 39.40000000000001
$ws.Cells.Item(41, 3).Value = 39.52373123168945
$ws.Cells.Item(41, 4).Value = 0.1237312316894474
$ws.Cells.Item(41, 5).Value = 0.01530941769538772
$ws.Cells.Item(42, 2).Value = 39.90000000000001
$ws.Cells.Item(42, 3).Value = 39.62896347045898
$ws.Cells.Item(42, 4).Value = -0.2710365295410213
$ws.Cells.Item(42, 5).Value = 0.07346080034564091
$ws.Cells.Item(43, 2).Value = 40.09999999999999
$ws.Cells.Item(43, 3).Value = 39.96549606323242
$ws.Cells.Item(43, 4).Value = -0.1345039367675724
$ws.Cells.Item(43, 5).Value = 0.01809130900597513
$ws.Cells.Item(44, 2).Value = 40.59999999999999
$ws.Cells.Item(44, 3).Value = 40.37540817260742
$ws.Cells.Item(44, 4).Value = -0.2245918273925724
$ws.Cells.Item(44, 5).Value = 0.05044148893153505
$ws.Cells.Item(45, 2).Value = 40.90000000000001
$ws.Cells.Item(45, 3).Value = 40.5744743347168
$ws.Cells.Item(45, 4).Value = -0.3255256652832088
$ws.Cells.Item(45, 5).Value = 0.1059669587580757
$ws.Cells.Item(46, 2).Value = 41.20000000000001
$ws.Cells.Item(46, 3).Value = 41.21307373046875
$ws.Cells.Item(46, 4).Value = 0.01307373046874005
$ws.Cells.Item(46, 5).Value = 0.000170922428369262
$ws.Cells.Item(47, 3).Value = 41.28318023681641
$ws.Cells.Item(47, 4).Value = -0.2168197631835938
$ws.Cells.Item(47, 5).Value = 0.04701080970698968
$ws.Cells.Item(48, 3).Value = 41.82052230834961
$ws.Cells.Item(48, 4).Value = 0.02052230834961222
$ws.Cells.Item(48, 5).Value = 0.0004211651399965633
$ws.Cells.Item(49, 3).Value = 41.57979202270508
$ws.Cells.Item(49, 4).Value = -0.6202079772949247
$ws.Cells.Item(49, 5).Value = 0.3846579351002619
$ws.Cells.Item(50, 3).Value = 43.23081970214844
$ws.Cells.Item(50, 4).Value = 0.5308197021484347
$ws.Cells.Item(50, 5).Value = 0.2817695561889529
$ws.Cells.Item(51, 3).Value = 43.8011360168457
$ws.Cells.Item(51, 4).Value = 0.1011360168457074
$ws.Cells.Item(51, 5).Value = 0.01022849390341521
$ws.Cells.Item(52, 3).Value = -0.07380088806155527
$ws.Cells.Item(52, 5).Value = 3.070135179046037
$ws.Cells.Item(53, 5).Value = 0.06140270358092074
